$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, pushing existing rows 32-131 down to 33-132.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new data point.
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44715
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112043
$ws.Range("G32").Value = "Pepino ensalada"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 140
$ws.Range("K32").Value = 18000
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 19143
$ws.Range("N32").Value = "$/caja 60 unidades"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 319
$ws.Range("Q32").Value = 60
$ws.Range("R32").Value = "Hortaliza"
